# Updates cryptos list (Price / Volume(1h) columns, and a swap of the
# TheGraph / LidoDAOToken rows) to match the "Updated cryptos list" commit.
# NumberFormat is forced to Text ("@") before assigning any cell whose new
# value would otherwise be auto-parsed by Excel as a number (which would
# silently drop things like trailing zeros, e.g. "10.50" -> 10.5).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.751.57'
$ws.Range('E2').Value = '  -0.60%  '
$ws.Range('D3').Value = '3.795.72'
$ws.Range('E3').Value = '  -1.30%  '
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '432.64'
$ws.Range('E5').Value = '  +4.87%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.07'
$ws.Range('E6').Value = '  +6.65%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.626'
$ws.Range('E7').Value = '  +1.16%  '
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.736'
$ws.Range('E9').Value = '  -0.51%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.153'
$ws.Range('E10').Value = '  -10.88%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000316'
$ws.Range('E11').Value = '  -16.97%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '42.95'
$ws.Range('E12').Value = '  +4.56%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.50'
$ws.Range('E13').Value = '  +4.28%  '
$ws.Range('D14').Value = '4.387.86'
$ws.Range('E14').Value = '  -1.46%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.98'
$ws.Range('E15').Value = '  -0.43%  '
$ws.Range('E16').Value = '  -0.18%  '
$ws.Range('D17').Value = '3.804.94'
$ws.Range('E17').Value = '  -1.05%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '19.98'
$ws.Range('E18').Value = '  +2.00%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.14'
$ws.Range('E19').Value = '  +5.58%  '
$ws.Range('D20').Value = '66.786.19'
$ws.Range('E20').Value = '  -1.07%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '410.49'
$ws.Range('E21').Value = '  -1.27%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '14.76'
$ws.Range('E22').Value = '  +0.67%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.26'
$ws.Range('E23').Value = '  +5.82%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '85.40'
$ws.Range('E24').Value = '  -0.79%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '36.86'
$ws.Range('E25').Value = '  -0.08%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.35'
$ws.Range('E26').Value = '  +6.75%  '
$ws.Range('E27').Value = '  -2.51%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.59'
$ws.Range('E28').Value = '  +32.46%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.81'
$ws.Range('E29').Value = '  +3.08%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '727.44'
$ws.Range('E30').Value = '  +6.43%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '13.85'
$ws.Range('E31').Value = '  +10.25%  '
$ws.Range('E32').Value = '  +9.95%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.70'
$ws.Range('E33').Value = '  -1.79%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '41.69'
$ws.Range('E34').Value = '  +6.80%  '
$ws.Range('E35').Value = '  +0.04%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.70'
$ws.Range('E36').Value = '  +27.09%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.151'
$ws.Range('E37').Value = '  -1.59%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '55.94'
$ws.Range('E38').Value = '  +1.23%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0476'
$ws.Range('E39').Value = '  +3.08%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.70'
$ws.Range('E40').Value = '  +39.36%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.93'
$ws.Range('E41').Value = '  -3.92%  '
$ws.Range('D42').Value = '0.0₃0682'
$ws.Range('E42').Value = '  -15.51%  '
$ws.Range('E43').Value = '  +2.78%  '
$ws.Range('E44').Value = '  +0.25%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.23'
$ws.Range('E45').Value = '  +1.40%  '
$ws.Range('B46').Value = 'LidoDAOToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.34'
$ws.Range('E46').Value = '  -0.12%  '
$ws.Range('B47').Value = 'TheGraph'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.320'
$ws.Range('E47').Value = '  +8.12%  '
$ws.Range('E48').Value = '  +3.76%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.10'
$ws.Range('E49').Value = '  -0.15%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '142.16'
$ws.Range('E50').Value = '  -4.52%  '
$ws.Range('E51').Value = '  -0.50%  '
